$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "In all, 43 risk estimates from 14 publications (2027 IPF cases in total) were used. Each exposure category was assessed with 6-11 risk estimates. Pooled ORs were significantly elevated for each category; the pooled PAF estimates by category ranged from 4-14% (Table 2). ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In all, 40 risk estimates from 12 publications (1326 IPF cases in total) were used. Each exposure category was assessed with 5-11 risk estimates. Pooled ORs were significantly elevated for each category; the pooled PAF estimates by category ranged from 3-23% (Table 2). ",
    2
)
